$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new innings/activity entry for Sanju Samson is recorded as the new row 7,
# pushing the previous last data row down to row 8.
$ws.Rows.Item(7).Insert()

# --- Populate the new row 7 with the freshly recorded activity -------------
# runs=0, balls=3, fours=0, sixes=0.
# Every numeric-looking cell in this sheet is stored as TEXT, not a number
# (see the original C2:F7 cells). Typing a numeric literal into .Value would
# make Excel auto-detect it as a Number, so instead we copy the values in
# from cells that already hold the same text (Range.Copy preserves the
# source cell's type), keeping the new cells text-typed just like the rest
# of the sheet.
$ws.Range("A6:B6").Copy($ws.Range("A7:B7"))   # player / team name
$ws.Range("F6").Copy($ws.Range("C7"))         # "0"
$ws.Range("F2").Copy($ws.Range("D7"))         # "3"
$ws.Range("E6").Copy($ws.Range("E7"))         # "0"
$ws.Range("F6").Copy($ws.Range("F7"))         # "0"

# --- Rows 5 and 8 trade places ---------------------------------------------
# After the insert, row 5 still has the old row-5 figures (9, 6, 0, 1) and
# row 8 has the old row-7 figures (74, 32, 1, 9) that got shifted down. The
# final sheet needs those two rows swapped, so stage row 5 off to one side,
# move row 8 into row 5, then move the staged values into row 8.
$stage = $ws.Range("Z1:AC1")
$ws.Range("C5:F5").Copy($stage)
$ws.Range("C8:F8").Copy($ws.Range("C5:F5"))
$stage.Copy($ws.Range("C8:F8"))
$stage.ClearContents()
